# Generate Report for handback
# Updates the zh-cn and de-de status sheets to reflect that both tracked
# files (0d410992-...md and 41cbefcb-...md) have been handed back and are
# in sync with en-US: fills in "Latest Target File" (E), "Latest Handback
# File" (F) and "Latest Handback DateTime" (G), and flips the Status (B)
# text from "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

# Cornflower blue (FF6495ED) packed as a BGR OLE color, matching the
# existing hyperlink-styled cells (A/C columns) in this workbook.
$linkColor = 15570276

function Apply-HandbackRow {
    param(
        $ws,
        $statusCell,
        $targetCell,
        $targetUrl,
        $targetDisplay,
        $handbackCell,
        $handbackUrl,
        $handbackDisplay,
        $datetimeCell,
        $datetimeValue
    )

    $ws.Range($statusCell).Value = "Handed back: in sync with en-US"

    $ws.Hyperlinks.Add($ws.Range($targetCell), $targetUrl, "", "", $targetDisplay)
    $ws.Range($targetCell).Font.Underline = 2
    $ws.Range($targetCell).Font.Color = $linkColor

    $ws.Hyperlinks.Add($ws.Range($handbackCell), $handbackUrl, "", "", $handbackDisplay)
    $ws.Range($handbackCell).Font.Underline = 2
    $ws.Range($handbackCell).Font.Color = $linkColor

    $ws.Range($datetimeCell).Value = $datetimeValue
    $ws.Range($datetimeCell).NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Apply-HandbackRow $wsZh `
    "B2" `
    "E2" `
    "https://github.com/OpenLocalizationTest/oltest/blob/1b739552f556e5a6f3111fd8148a8f0ff47abede/e2e/0d410992-14c1-46f6-b249-15f558031b5a.md" `
    "0d410992-14c1-46f6-b249-15f558031b5a.md" `
    "F2" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d86535715fd5fc49d9462ae36b67d332e4408e32/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/0d410992-14c1-46f6-b249-15f558031b5a.d6ebf941e230af1fb2b251765401ee991c989ffb.zh-cn.xlf" `
    "0d410992-14c1-46f6-b249-15f558031b5a.d6ebf941e230af1fb2b251765401ee991c989ffb.zh-cn.xlf" `
    "G2" `
    "2016-01-13 04:44:00"

Apply-HandbackRow $wsZh `
    "B3" `
    "E3" `
    "https://github.com/OpenLocalizationTest/oltest/blob/1b739552f556e5a6f3111fd8148a8f0ff47abede/e2e/41cbefcb-0561-4345-abaf-f0b90f69efd9.md" `
    "41cbefcb-0561-4345-abaf-f0b90f69efd9.md" `
    "F3" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d86535715fd5fc49d9462ae36b67d332e4408e32/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/41cbefcb-0561-4345-abaf-f0b90f69efd9.3673da69609661c722dc95eb85a098e28a3e6cfd.zh-cn.xlf" `
    "41cbefcb-0561-4345-abaf-f0b90f69efd9.3673da69609661c722dc95eb85a098e28a3e6cfd.zh-cn.xlf" `
    "G3" `
    "2016-01-13 04:44:00"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Apply-HandbackRow $wsDe `
    "B2" `
    "E2" `
    "https://github.com/OpenLocalizationTest/oltest/blob/1b739552f556e5a6f3111fd8148a8f0ff47abede/e2e/0d410992-14c1-46f6-b249-15f558031b5a.md" `
    "0d410992-14c1-46f6-b249-15f558031b5a.md" `
    "F2" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8949e8f30d431eda4b0d4d067b868032753f3105/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/0d410992-14c1-46f6-b249-15f558031b5a.d6ebf941e230af1fb2b251765401ee991c989ffb.de-de.xlf" `
    "0d410992-14c1-46f6-b249-15f558031b5a.d6ebf941e230af1fb2b251765401ee991c989ffb.de-de.xlf" `
    "G2" `
    "2016-01-13 04:44:34"

Apply-HandbackRow $wsDe `
    "B3" `
    "E3" `
    "https://github.com/OpenLocalizationTest/oltest/blob/1b739552f556e5a6f3111fd8148a8f0ff47abede/e2e/41cbefcb-0561-4345-abaf-f0b90f69efd9.md" `
    "41cbefcb-0561-4345-abaf-f0b90f69efd9.md" `
    "F3" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8949e8f30d431eda4b0d4d067b868032753f3105/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/41cbefcb-0561-4345-abaf-f0b90f69efd9.3673da69609661c722dc95eb85a098e28a3e6cfd.de-de.xlf" `
    "41cbefcb-0561-4345-abaf-f0b90f69efd9.3673da69609661c722dc95eb85a098e28a3e6cfd.de-de.xlf" `
    "G3" `
    "2016-01-13 04:44:34"
